$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ыы"
$ws.Range("B2").Value = "2 шт"
$ws.Range("C2").Value = "2 шт"
$ws.Range("D2").Value = 40.0

$ws.Range("A3").Value = "ИТОГО:"
$ws.Range("D3").Value = 40.0
